$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append, starting at row 256 (date serial 44330, 2021-05-14)
# through row 269 (date serial 44343, 2021-05-27), matching the commit
# "aggiornamento fino a 27/05".
$newData = @(
    @(44330, 0, 0, 0),
    @(44331, 1, 1, 43.78283712784589),
    @(44332, 1, 2, 87.56567425569177),
    @(44333, 0, 2, 87.56567425569177),
    @(44334, 2, 4, 175.1313485113835),
    @(44335, 0, 4, 175.1313485113835),
    @(44336, 0, 4, 175.1313485113835),
    @(44337, 0, 4, 175.1313485113835),
    @(44338, 0, 3, 131.3485113835376),
    @(44339, 0, 2, 87.56567425569177),
    @(44340, 0, 2, 87.56567425569177),
    @(44341, 0, 0, 0),
    @(44342, 0, 0, 0),
    @(44343, 0, 0, 0)
)

$startRow = 256
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]

    $srcCell = $ws.Cells.Item($r - 1, 1)
    $cellA = $ws.Cells.Item($r, 1)
    $srcCell.Copy($cellA)
    $cellA.Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
